# Add Giftcode in My_Account Osprey
# On the "Giftcard Payments" sheet: insert a new column before column J,
# then populate the new header (J1) and the new row-22 value (J22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Giftcard Payments")

# Insert a new blank column before J (shifts J:AK -> K:AL)
$ws.Columns("J:J").Insert()

# New header cell for the inserted column
$ws.Range("J1").Value = "GiftCard2"

# New value cell in row 22 for the inserted column
$ws.Range("J22").Value = "OSPREY-GIFT-CARD_ST827"

# Leave the new column's header cell selected (matches the saved selection state)
$ws.Range("J1").Select() | Out-Null
